$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 144
$ws.Range("I2").Value = 335
$ws.Range("J2").Value = 1532
$ws.Range("L2").Value = 456
$ws.Range("M2").Value = 17
$ws.Range("N2").Value = 263
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 19
$ws.Range("S2").Value = 190
$ws.Range("T2").Value = 277
$ws.Range("U2").Value = 19
$ws.Range("V2").Value = 2521
$ws.Range("X2").Value = 2416
$ws.Range("Y2").Value = 4
$ws.Range("Z2").Value = 37
$ws.Range("AA2").Value = 19
